# SWOT Matrix: split the single 2x2 SWOT slide into two slides -
#   Slide 1 ("Internal Factors"): Strengths + Weaknesses, enlarged to
#     fill the vertical space previously shared with Opportunities/Threats.
#   Slide 2 ("External Factors", new): Opportunities + Threats, moved up
#     into the (now vacated) top row that Strengths/Weaknesses occupy on
#     slide 1.
#
# NOTE on units: Shape.Left/Top/Width/Height (and AddTextbox's geometry
# args) are expressed in points (1 pt = 12700 EMU) in this object model.
# The runtime stores these coordinates with float32 precision internally,
# so a plain `emu / 12700.0` assignment can truncate to one EMU below the
# intended value (e.g. 365760 EMU becomes 365759 EMU) when the value is
# fed back through the Left/Top/Width/Height *setters* (this does not
# happen for the literal geometry passed straight into AddTextbox). Adding
# a hair of an EMU before converting keeps the float32 round-trip on the
# correct side for every offset/extent this script needs, without ever
# being large enough to push a value into the next EMU.
function EmuToPt {
    param([double]$Emu)
    return ($Emu / 12700.0) + 0.00003937007874015748
}

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Slide 1: retitle, and grow the Strengths/Weaknesses content boxes down
# into the space that Opportunities/Threats used to occupy.
# ---------------------------------------------------------------------
$titleShape = $s1.Shapes.Item(1)
$titleRun = $titleShape.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$titleRun.Text = "SWOT Analysis " + [char]0x2014 + " Internal Factors"
# Re-assert the original box height; setting the run text re-triggers the
# shape's spAutoFit recalculation, which isn't what the source edit did.
$titleShape.Height = EmuToPt(548640)

$strengthsContent = $s1.Shapes.Item(3)
$strengthsContent.Height = EmuToPt(5093208)

$weaknessesContent = $s1.Shapes.Item(5)
$weaknessesContent.Height = EmuToPt(5093208)

# ---------------------------------------------------------------------
# Pull off the Opportunities/Threats shapes (still on slide 1, at indices
# 6-9) before they're deleted, so slide 2 can be built from exact copies
# that keep every bit of their original formatting.
# ---------------------------------------------------------------------
$oppHeaderSrc = $s1.Shapes.Item(6)
$oppContentSrc = $s1.Shapes.Item(7)
$threatHeaderSrc = $s1.Shapes.Item(8)
$threatContentSrc = $s1.Shapes.Item(9)

# New slide, inserted right after slide 1, on the Blank layout (matching
# slide 1's own layout, so no placeholder shapes come along for free).
$s2 = $p.Slides.Add(2, 12)

# Title + subtitle textbox, copied from slide 1's title so formatting
# (fonts, fills, auto-fit) matches exactly; text is then swapped.
$titleShape.Copy()
$newTitle = $s2.Shapes.Paste().Item(1)
$newTitle.Name = "TextBox 1"
$newTitleRun = $newTitle.TextFrame.TextRange.Paragraphs(1, 1).Runs(1, 1)
$newTitleRun.Text = "SWOT Analysis " + [char]0x2014 + " External Factors"
$newTitle.Left = EmuToPt(365760)
$newTitle.Top = EmuToPt(228600)
$newTitle.Width = EmuToPt(11430000)
$newTitle.Height = EmuToPt(548640)

# Opportunities header -> slide 2, top-left (where Strengths' header was).
$oppHeaderSrc.Copy()
$newOppHeader = $s2.Shapes.Paste().Item(1)
$newOppHeader.Name = "TextBox 2"
$newOppHeader.Left = EmuToPt(365760)
$newOppHeader.Top = EmuToPt(1051560)

# Opportunities content -> slide 2, enlarged the same way the internal
# factors boxes were on slide 1.
$oppContentSrc.Copy()
$newOppContent = $s2.Shapes.Paste().Item(1)
$newOppContent.Name = "TextBox 3"
$newOppContent.Left = EmuToPt(365760)
$newOppContent.Top = EmuToPt(1399031)
$newOppContent.Width = EmuToPt(5577840)
$newOppContent.Height = EmuToPt(5093208)

# Threats header -> slide 2, top-right (where Weaknesses' header was).
$threatHeaderSrc.Copy()
$newThreatHeader = $s2.Shapes.Paste().Item(1)
$newThreatHeader.Name = "TextBox 4"
$newThreatHeader.Left = EmuToPt(6309360)
$newThreatHeader.Top = EmuToPt(1051560)

# Threats content -> slide 2, enlarged to match.
$threatContentSrc.Copy()
$newThreatContent = $s2.Shapes.Paste().Item(1)
$newThreatContent.Name = "TextBox 5"
$newThreatContent.Left = EmuToPt(6309360)
$newThreatContent.Top = EmuToPt(1399031)
$newThreatContent.Width = EmuToPt(5577840)
$newThreatContent.Height = EmuToPt(5093208)

# Now that slide 2 has its own copies, drop the originals from slide 1
# (delete back-to-front so earlier indices stay valid).
$threatContentSrc.Delete()
$threatHeaderSrc.Delete()
$oppContentSrc.Delete()
$oppHeaderSrc.Delete()
